# Implement ICP + LRP - first version
# Insert a new task row ("Write LRP algorithm") above the "DI" row in the
# weekly point table (around row 120), shifting the rest of the table down
# by one row, fill in the new row's data, add the missing "MO" points value,
# and keep the conditional-formatting / selection state in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# --- Insert a new row at 120 (existing rows 120:135 shift down to 121:136) ---
$ws.Rows.Item(120).Insert()

# --- Fill the new row 120 with the new task ---
$ws.Range("E120").Value = "Write LRP algorithm"
$ws.Range("F120").Value = 6
$ws.Range("E120:F120").Interior.Color = 65535

# --- Row 118 (E118:F118) picks up the same highlight color as E117:F117 ---
$ws.Range("E118:F118").Interior.Color = 1048374

# --- Fill in the previously-missing "MO" point value (now row 125) ---
$ws.Range("J125").Value = 8

# --- Conditional formatting ranges follow the rows they were anchored to ---
$ws.Range("J127").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J128"))
$ws.Range("J135").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J136"))

# --- Update view state: scrolled position / active selection ---
$ws.Range("K130").Select() | Out-Null
